# Final Project Submission in Fulfillment of Requirements
#
# The author re-saved the workbook from a local (non-SharePoint-synced)
# copy of Excel: the sheet "final model reduced features" was the active
# sheet, and its previous block selection (B16:G25) was collapsed down to
# a single-cell selection on A16 before the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("final model reduced features")

# Keep/confirm this sheet as the active one (it already was).
$ws.Activate()

# Collapse the old B16:G25 block selection down to a single active cell, A16.
$ws.Range("A16").Select()
